# edit.ps1 - apply the "Colors.xlsx" commit:
#  - fill in previously-blank H-column cells for rows 50,51,52,53,54,58
#  - expand the text of H56 and H61
#  - change the active selection / scroll position of the sheet view
#    (from topLeftCell A46 / selection H56  ->  no frozen topLeftCell / selection G17)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- New values for previously empty H cells -----------------------------
$ws.Range("H50").Value = "두록색;회색;백색;자황색;적색;청색;흑청색"
$ws.Range("H51").Value = "백색;치색;벽자색;청색;군청색;명황색;적색"
$ws.Range("H52").Value = "소색;치색;백색;청벽색;벽청색;흑색"
$ws.Range("H53").Value = "구색;흑색;백색;연지회색;청색;소색"
$ws.Range("H54").Value = "소색;다자색;회색;담자색;적색;구색;흑색"
$ws.Range("H58").Value = "지백색;자황색;갈색;백색;연지회색;유황색"

# --- Expanded text for existing H cells -----------------------------------
$ws.Range("H56").Value = "백색;자황색;소색;홍황색;구색;연지회색;자황색;연두색"
$ws.Range("H61").Value = "황색;벽청색;비색;청색;유록색;청자색"

# --- Update sheet view: scroll back to top-left and move selection --------
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("G17").Select()
